$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("PWM")
$ws2 = $wb.Worksheets.Item("ADC")

# ---------------------------------------------------------------
# Sheet1 (PWM): insert a new "Resolution, bit" row above the old
# "TOP" row, turning TOP into a formula, and append a new duty
# cycle calculation block in rows 8-11.
# ---------------------------------------------------------------

# Push the existing TOP/Fpwm rows down by inserting a blank row at 3
$ws1.Rows("3").Insert()

$ws1.Range("A3").Value = "Resolution, bit"
$ws1.Range("B3").Value = 10

# Old row 3 (TOP) is now row 4; turn its literal value into a formula
$ws1.Range("B4").Formula = "=POWER(2,10)-1"

# Old row 4 (Fpwm) is now row 5; Insert() already shifted its formula's
# cell reference from B3 to B4 automatically, nothing further needed.

# New duty-cycle block
$ws1.Range("A8").Value = "Ubat, V"
$ws1.Range("B8").Value = 13.6

$ws1.Range("A9").Value = "condens, V"
$ws1.Range("B9").Value = 1.5

$ws1.Range("A10").Value = "duty, %"
$ws1.Range("B10").Formula = "=(B9/B8)*100"

$ws1.Range("A11").Value = "dutyCycle"
$ws1.Range("B11").Formula = "=(1024*B10)/100"

# Column A width on PWM sheet
$ws1.Columns("A").ColumnWidth = 13.14

# ---------------------------------------------------------------
# Sheet2 (ADC): corrected ADC raw reading
# ---------------------------------------------------------------
$ws2.Range("B2").Value = 563

# ---------------------------------------------------------------
# View / selection state: ADC becomes the active/selected tab,
# PWM keeps a different selected cell.
# ---------------------------------------------------------------
$ws1.Range("A17").Select()
$ws2.Activate()
$ws2.Range("B3").Select()
